# Add signal for the out door
# Row 67 (terminal № 64) gets the new "Out Door Tamper" signal entry,
# and the remaining spare rows (68-78) get their Terminal column
# pre-formatted as text, matching the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new row (and the remaining spare Terminal cells) as text,
# same "Text" number format used throughout the table, before typing
# the values so "12.1" is stored as text, not a number.
$ws.Range("B67:F67").NumberFormat = "@"
$ws.Range("B68:B78").NumberFormat = "@"

# New signal row: Terminal, Signal Type, Signal Name, Function, Location
# (written Signal Name / Function first, then Terminal, matching the
# order new shared strings were registered in the source file)
$ws.Range("D67").Value = "ODT"
$ws.Range("E67").Value = "Out Door Tamper"
$ws.Range("B67").Value = "12.1"
$ws.Range("C67").Value = "DI-0-24V"
$ws.Range("F67").Value = "Cabinet"
